# Auto-generated Excel COM-interop script applying the Sargatanas_Profits value refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 654.2941
$ws.Cells.Item(17, 10).Value = 748.14636
$ws.Cells.Item(17, 12).Value = 2244.43908
$ws.Cells.Item(17, 14).Value = -2580.43908
$ws.Cells.Item(40, 8).Value = 2699.7
$ws.Cells.Item(40, 9).Value = 3083.3333
$ws.Cells.Item(40, 10).Value = 2124.25
$ws.Cells.Item(40, 11).Value = 3083.3333
$ws.Cells.Item(40, 12).Value = 2124.25
$ws.Cells.Item(40, 13).Value = -2908.3333
$ws.Cells.Item(40, 14).Value = -2474.25
$ws.Cells.Item(70, 8).Value = 216668180
$ws.Cells.Item(70, 9).Value = 83335830
$ws.Cells.Item(70, 11).Value = 250007490
$ws.Cells.Item(70, 13).Value = -250007220
$ws.Cells.Item(73, 8).Value = 216668180
$ws.Cells.Item(73, 9).Value = 83335830
$ws.Cells.Item(73, 11).Value = 250007490
$ws.Cells.Item(73, 13).Value = -250006554
$ws.Cells.Item(76, 8).Value = 203999.2
$ws.Cells.Item(76, 9).Value = 502499
$ws.Cells.Item(76, 10).Value = 4999.3335
$ws.Cells.Item(76, 11).Value = 502499
$ws.Cells.Item(76, 12).Value = 4999.3335
$ws.Cells.Item(76, 13).Value = -502184
$ws.Cells.Item(76, 14).Value = -5629.3335
$ws.Cells.Item(79, 8).Value = 203999.2
$ws.Cells.Item(79, 9).Value = 502499
$ws.Cells.Item(79, 10).Value = 4999.3335
$ws.Cells.Item(79, 11).Value = 502499
$ws.Cells.Item(79, 12).Value = 4999.3335
$ws.Cells.Item(79, 13).Value = -501407
$ws.Cells.Item(79, 14).Value = -7183.3335
$ws.Cells.Item(80, 8).Value = 8367058
$ws.Cells.Item(80, 10).Value = 100625
$ws.Cells.Item(80, 12).Value = 301875
$ws.Cells.Item(80, 14).Value = -303871
$ws.Cells.Item(83, 8).Value = 8367058
$ws.Cells.Item(83, 10).Value = 100625
$ws.Cells.Item(83, 12).Value = 905625
$ws.Cells.Item(83, 14).Value = -915609
$ws.Cells.Item(86, 8).Value = 93752570
$ws.Cells.Item(86, 9).Value = 111113000
$ws.Cells.Item(86, 10).Value = 18524024
$ws.Cells.Item(86, 11).Value = 111113000
$ws.Cells.Item(86, 12).Value = 18524024
$ws.Cells.Item(86, 13).Value = -111111877
$ws.Cells.Item(86, 14).Value = -18526270
$ws.Cells.Item(89, 8).Value = 93752570
$ws.Cells.Item(89, 9).Value = 111113000
$ws.Cells.Item(89, 10).Value = 18524024
$ws.Cells.Item(89, 11).Value = 555565000
$ws.Cells.Item(89, 12).Value = 92620120
$ws.Cells.Item(89, 13).Value = -555559384
$ws.Cells.Item(89, 14).Value = -92631352
$ws.Cells.Item(107, 8).Value = 19445526
$ws.Cells.Item(107, 10).Value = 58334784
$ws.Cells.Item(107, 12).Value = 58334784
$ws.Cells.Item(107, 14).Value = -58338624
$ws.Cells.Item(137, 8).Value = 4465.4443
$ws.Cells.Item(137, 9).Value = 20000
$ws.Cells.Item(137, 10).Value = 2523.625
$ws.Cells.Item(137, 11).Value = 60000
$ws.Cells.Item(137, 12).Value = 7570.875
$ws.Cells.Item(137, 13).Value = -57450
$ws.Cells.Item(137, 14).Value = -12670.875
$ws.Cells.Item(138, 8).Value = 5462.108
$ws.Cells.Item(138, 9).Value = 1595.4
$ws.Cells.Item(138, 10).Value = 10011.177
$ws.Cells.Item(138, 11).Value = 4786.200000000001
$ws.Cells.Item(138, 12).Value = 30033.531
$ws.Cells.Item(138, 13).Value = 353.7999999999993
$ws.Cells.Item(138, 14).Value = -40313.531

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 33032.434
$ws.Cells.Item(74, 9).Value = 45807.863
$ws.Cells.Item(74, 11).Value = 45807.863
$ws.Cells.Item(74, 13).Value = -44933.863
$ws.Cells.Item(77, 8).Value = 33032.434
$ws.Cells.Item(77, 9).Value = 45807.863
$ws.Cells.Item(77, 11).Value = 229039.315
$ws.Cells.Item(77, 13).Value = -224671.315
$ws.Cells.Item(132, 8).Value = 4708.254
$ws.Cells.Item(132, 9).Value = 3464.9092
$ws.Cells.Item(132, 11).Value = 10394.7276
$ws.Cells.Item(132, 13).Value = -7864.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 10418376
$ws.Cells.Item(20, 9).Value = 12822308
$ws.Cells.Item(20, 10).Value = 1335.3334
$ws.Cells.Item(20, 11).Value = 12822308
$ws.Cells.Item(20, 12).Value = 1335.3334
$ws.Cells.Item(20, 13).Value = -12822061
$ws.Cells.Item(20, 14).Value = -1829.3334
$ws.Cells.Item(134, 8).Value = 5001.3584
$ws.Cells.Item(134, 9).Value = 1890.875
$ws.Cells.Item(134, 11).Value = 5672.625
$ws.Cells.Item(134, 13).Value = -3137.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9624.611000000001
$ws.Cells.Item(31, 9).Value = 3733.6428
$ws.Cells.Item(31, 10).Value = 13373.409
$ws.Cells.Item(31, 11).Value = 3733.6428
$ws.Cells.Item(31, 12).Value = 13373.409
$ws.Cells.Item(31, 13).Value = -3438.6428
$ws.Cells.Item(31, 14).Value = -13963.409
$ws.Cells.Item(34, 8).Value = 9624.611000000001
$ws.Cells.Item(34, 9).Value = 3733.6428
$ws.Cells.Item(34, 10).Value = 13373.409
$ws.Cells.Item(34, 11).Value = 3733.6428
$ws.Cells.Item(34, 12).Value = 13373.409
$ws.Cells.Item(34, 13).Value = -3531.6428
$ws.Cells.Item(34, 14).Value = -13777.409
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 14).Value = 0
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 14).Value = 0
$ws.Cells.Item(99, 8).Value = 11083.25
$ws.Cells.Item(99, 10).Value = 7875
$ws.Cells.Item(99, 12).Value = 7875
$ws.Cells.Item(99, 14).Value = -10871
$ws.Cells.Item(126, 8).Value = 11083.25
$ws.Cells.Item(126, 10).Value = 7875
$ws.Cells.Item(126, 12).Value = 23625
$ws.Cells.Item(126, 14).Value = -28565
$ws.Cells.Item(132, 8).Value = 4963.0933
$ws.Cells.Item(132, 9).Value = 2258.6667
$ws.Cells.Item(132, 11).Value = 6776.000100000001
$ws.Cells.Item(132, 13).Value = -4246.000100000001
$ws.Cells.Item(139, 8).Value = 65499.5
$ws.Cells.Item(139, 10).Value = 69599.39999999999
$ws.Cells.Item(139, 12).Value = 69599.39999999999
$ws.Cells.Item(139, 14).Value = -79879.39999999999
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 5863.25
$ws.Cells.Item(56, 9).Value = 5863.25
$ws.Cells.Item(56, 11).Value = 5863.25
$ws.Cells.Item(56, 13).Value = -5333.25
$ws.Cells.Item(129, 8).Value = 83833336
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 83833336
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 13).Value = 251500008
$ws.Cells.Item(129, 14).Value = -251510008
$ws.Cells.Item(141, 8).Value = 7019.1333
$ws.Cells.Item(141, 9).Value = 2587.5557
$ws.Cells.Item(141, 11).Value = 7762.6671
$ws.Cells.Item(141, 13).Value = -2582.6671
$ws.Range("M129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 5393.696
$ws.Cells.Item(102, 9).Value = 4370
$ws.Cells.Item(102, 10).Value = 9079
$ws.Cells.Item(102, 11).Value = 4370
$ws.Cells.Item(102, 12).Value = 9079
$ws.Cells.Item(102, 13).Value = -2748
$ws.Cells.Item(102, 14).Value = -12323

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7213.7856
$ws.Cells.Item(40, 9).Value = 4997.25
$ws.Cells.Item(40, 11).Value = 4997.25
$ws.Cells.Item(40, 13).Value = -4861.25
$ws.Cells.Item(132, 8).Value = 10007211
$ws.Cells.Item(132, 9).Value = 16671167
$ws.Cells.Item(132, 11).Value = 50013501
$ws.Cells.Item(132, 13).Value = -50010971

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1646.7
$ws.Cells.Item(96, 9).Value = 1682.1666
$ws.Cells.Item(96, 10).Value = 1593.5
$ws.Cells.Item(96, 11).Value = 1682.1666
$ws.Cells.Item(96, 12).Value = 1593.5
$ws.Cells.Item(96, 13).Value = -309.1666
$ws.Cells.Item(96, 14).Value = -4339.5
$ws.Cells.Item(132, 8).Value = 20851776
$ws.Cells.Item(132, 9).Value = 38471604
$ws.Cells.Item(132, 10).Value = 28341.545
$ws.Cells.Item(132, 11).Value = 115414812
$ws.Cells.Item(132, 12).Value = 85024.63499999999
$ws.Cells.Item(132, 13).Value = -115412282
$ws.Cells.Item(132, 14).Value = -90084.63499999999
